$wb = $excel.ActiveWorkbook

# --- Sheet3 -> renamed to "testDemo", populate with User/Pass data ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Name = "testDemo"
$ws3.Activate()
$ws3.Range("A1").Value = "User"
$ws3.Range("B1").Value = "Pass"
$ws3.Range("A2").Value = 454
$ws3.Range("B2").Value = 554
$ws3.Range("A3").Value = 454
$ws3.Range("B3").Value = 554
[void]$ws3.Range("B3").Select()

# --- Sheet2 (validLoginTest) -> populate with login data + hyperlink, becomes active tab ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Activate()
$ws2.Range("A1").Value = "Username"
$ws2.Range("B1").Value = "Password"
$ws2.Range("C1").Value = "ExpectedUrl"
$ws2.Hyperlinks.Add($ws2.Range("C2"), "https://opensource-demo.orangehrmlive.com/web/index.php/pim/viewEmployeeList")
$ws2.Range("A2").Value = "Admin"
$ws2.Range("B2").Value = "admin123"
[void]$ws2.Range("B3").Select()
